# Refresh the cryptos list (price + 1h volume change columns) to match
# the latest GitHub Actions scrape. Coin rows 40/41 and 44/45 also swap
# rank order (Celestia <-> Stacks, EnergySwap <-> TheGraph).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price-column (D) cells that look like plain numbers to stay text,
# matching the source sheet where every Price cell is a literal string
# (e.g. thousand-dot formatted "50.973.92", or trailing-zero "0.170").
$ws.Range("D2").Value = "50.973.92"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "2.943.76"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "375.81"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.72"
$ws.Range("E6").Value = "  -2.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.535"
$ws.Range("E7").Value = "  -0.92%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.586"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.31"
$ws.Range("E10").Value = "  -2.29%  "
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0839"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "3.402.64"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.94"
$ws.Range("E14").Value = "  -2.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.45"
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").Value = "2.936.03"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.981"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.81"
$ws.Range("E18").Value = "  +46.73%  "
$ws.Range("D19").Value = "50.874.71"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.15"
$ws.Range("E20").Value = "  -5.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.62"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").Value = "0.0₃0954"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "264.32"
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.48"
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.14"
$ws.Range("E25").Value = "  +10.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.25"
$ws.Range("E26").Value = "  +7.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.84"
$ws.Range("E27").Value = "  +6.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.170"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.64"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.89"
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.72"
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.58"
$ws.Range("E34").Value = "  -3.51%  "
$ws.Range("E35").Value = "  -2.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0442"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.07"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.116"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.39"
$ws.Range("E40").Value = "  -4.93%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.52"
$ws.Range("E41").Value = "  -2.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.79"
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.49"
$ws.Range("E43").Value = "  -2.69%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.19"
$ws.Range("E44").Value = "  -3.98%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.285"
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("E46").Value = "  -2.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.30"
$ws.Range("E47").Value = "  +2.62%  "
$ws.Range("E48").Value = "  -3.52%  "
$ws.Range("D49").Value = "1.994.42"
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0340"
$ws.Range("E50").Value = "  -2.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.28"
$ws.Range("E51").Value = "  -0.65%  "
